$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Contact Number column holds values with leading zeros (e.g.
# "03303184100"). Excel's Range.Value setter auto-detects numeric-looking
# strings and would otherwise strip the leading zero, so mark the column
# as Text first to preserve the values exactly as authored.
$ws.Range("E2:E4").NumberFormat = "@"

$ws.Range("C2").Value = "AUTODSR_735C4"
$ws.Range("D2").Value = "Father_557C"
$ws.Range("E2").Value = "03303184100"
$ws.Range("G2").Value = "EMP5A21C3"

$ws.Range("C3").Value = "AUTODSR_68BEB"
$ws.Range("D3").Value = "Father_9B3A"
$ws.Range("E3").Value = "03305402500"
$ws.Range("G3").Value = "EMP5C5B75"

$ws.Range("C4").Value = "AUTODSR_1204D"
$ws.Range("D4").Value = "Father_A820"
$ws.Range("E4").Value = "03307303800"
$ws.Range("G4").Value = "EMP211591"
